# Developer Guide: Fix spelling error
#  - Correct "CrearCommand" -> "ClearCommand" on the Design diagram slide.
#  - Refresh the cached "Fixed" date shown in the slide master / layout
#    footers (8/7/2018 -> 9/21/18).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer "date" placeholder text, on the slide master and every
#    custom (slide) layout that carries one.
# ---------------------------------------------------------------------
$oldDate = "8/7/2018"
$newDate = "9/21/18"

function Update-DateFooter($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateFooter $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateFooter $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Spelling fix: "CrearCommand" -> "ClearCommand" in the class diagram
#    on slide 2.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(2)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "CrearCommand") {
            $shp.TextFrame.TextRange.Text = "ClearCommand"
        }
    }
}
